$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must stay TEXT (as in the source diff).
# Temporarily mark them as Text-formatted before assigning, then restore the default "Normal"
# style so the saved XML does not pick up a stray style index on these cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.926.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.894.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7707"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3123"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.63"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07328"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.477"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.882.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.210"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.934.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007807"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.151.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.115"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1574"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.442"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.026"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.478"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05549"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.057"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.238"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7517"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.684"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01928"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.797"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4463"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.101.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.972"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8504"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.884"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.512"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.751"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.001"
$ws.Range("D51").Style = "Normal"

# Columns B, C and E are plain text (names, links, percentages) and do not need protection.
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("E5").Value = "  -3.48%  "
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("E10").Value = "  +4.02%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("E13").Value = "  +3.41%  "
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("E16").Value = "  +4.67%  "
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("E25").Value = "  -5.71%  "
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("E30").Value = "  +2.84%  "
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("E43").Value = "  +7.04%  "
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("E51").Value = "  +2.47%  "
